# Update crypto price/volume figures per the latest scrape.
# Values are written as literal text (via a scratch formula + paste-values)
# so Excel does not auto-convert numeric-looking strings (e.g. "0.999",
# "18.10", "58.267.38") into actual numbers, which would silently drop
# formatting (trailing zeros) and change the cell type / style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

$scratch.Formula = "=""58.267.38"""
$scratch.Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.53%  """
$scratch.Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""2.453.09"""
$scratch.Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.68%  """
$scratch.Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -0.15%  """
$scratch.Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""526.76"""
$scratch.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -2.83%  """
$scratch.Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""134.08"""
$scratch.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -8.39%  """
$scratch.Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.999"""
$scratch.Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  +0.50%  """
$scratch.Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.552"""
$scratch.Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.74%  """
$scratch.Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""2.458.88"""
$scratch.Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -4.81%  """
$scratch.Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.0989"""
$scratch.Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.16%  """
$scratch.Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -0.39%  """
$scratch.Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""5.31"""
$scratch.Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.47%  """
$scratch.Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.342"""
$scratch.Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -6.28%  """
$scratch.Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""2.887.65"""
$scratch.Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.63%  """
$scratch.Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""58.147.39"""
$scratch.Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.65%  """
$scratch.Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""22.47"""
$scratch.Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -8.23%  """
$scratch.Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -4.50%  """
$scratch.Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""2.457.38"""
$scratch.Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.07%  """
$scratch.Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""10.61"""
$scratch.Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -6.60%  """
$scratch.Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""318.53"""
$scratch.Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.19%  """
$scratch.Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""4.15"""
$scratch.Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -5.22%  """
$scratch.Copy() | Out-Null
$ws.Range("E21").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.996"""
$scratch.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -0.26%  """
$scratch.Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""5.67"""
$scratch.Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -5.30%  """
$scratch.Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""62.07"""
$scratch.Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -1.75%  """
$scratch.Copy() | Out-Null
$ws.Range("E24").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.405"""
$scratch.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -8.14%  """
$scratch.Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -3.10%  """
$scratch.Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.981"""
$scratch.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -0.98%  """
$scratch.Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""7.46"""
$scratch.Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -7.81%  """
$scratch.Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.0₃0747"""
$scratch.Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -8.10%  """
$scratch.Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""6.45"""
$scratch.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -9.84%  """
$scratch.Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""1.74"""
$scratch.Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -4.57%  """
$scratch.Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""163.11"""
$scratch.Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -1.25%  """
$scratch.Copy() | Out-Null
$ws.Range("E32").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.998"""
$scratch.Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  +0.13%  """
$scratch.Copy() | Out-Null
$ws.Range("E33").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -13.18%  """
$scratch.Copy() | Out-Null
$ws.Range("E34").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""18.10"""
$scratch.Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -4.10%  """
$scratch.Copy() | Out-Null
$ws.Range("E35").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -9.75%  """
$scratch.Copy() | Out-Null
$ws.Range("E36").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""3.98"""
$scratch.Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -11.45%  """
$scratch.Copy() | Out-Null
$ws.Range("E37").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""1.53"""
$scratch.Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -7.19%  """
$scratch.Copy() | Out-Null
$ws.Range("E38").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""36.38"""
$scratch.Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -1.92%  """
$scratch.Copy() | Out-Null
$ws.Range("E39").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""3.50"""
$scratch.Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -7.22%  """
$scratch.Copy() | Out-Null
$ws.Range("E40").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.776"""
$scratch.Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -8.22%  """
$scratch.Copy() | Out-Null
$ws.Range("E41").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.998"""
$scratch.Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  +0.64%  """
$scratch.Copy() | Out-Null
$ws.Range("E42").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""269.47"""
$scratch.Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -11.97%  """
$scratch.Copy() | Out-Null
$ws.Range("E43").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""4.98"""
$scratch.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -12.89%  """
$scratch.Copy() | Out-Null
$ws.Range("E44").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  +0.02%  """
$scratch.Copy() | Out-Null
$ws.Range("E45").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.585"""
$scratch.Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -4.40%  """
$scratch.Copy() | Out-Null
$ws.Range("E46").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.0917"""
$scratch.Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -2.71%  """
$scratch.Copy() | Out-Null
$ws.Range("E47").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""120.22"""
$scratch.Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -4.95%  """
$scratch.Copy() | Out-Null
$ws.Range("E48").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""0.0500"""
$scratch.Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -4.97%  """
$scratch.Copy() | Out-Null
$ws.Range("E49").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -6.41%  """
$scratch.Copy() | Out-Null
$ws.Range("E50").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""16.85"""
$scratch.Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$scratch.Formula = "=""  -8.90%  """
$scratch.Copy() | Out-Null
$ws.Range("E51").PasteSpecial(-4163) | Out-Null

$scratch.ClearContents() | Out-Null
$excel.CutCopyMode = 0

